$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","I","J","K","M","N")

$data = @(
  @(0.6743817574854631, 0.2449791818675635, 0.05732154878467099, 0.107454611188782, 2.878803597020351, 0.07973214163530429, 1.501337203754353, 0.1807450889423592, 0.9121141558725867, 0.3575122264537711, 2.967849150855301),
  @(0.6406894837234347, 0.2365572323609229, 0.0563272542536879, 0.1067583391954265, 2.871212034085957, 0.07973214163530429, 1.501496086960103, 0.180229736806524, 0.8708599391782172, 0.3475142182274098, 2.984003695558677),
  @(0.6203843842382923, 0.2315308718672782, 0.05574239339158993, 0.1063766714217742, 2.867788914015279, 0.07973214163530429, 1.502191119430307, 0.1799882657656084, 0.8460555906857792, 0.341572948509139, 2.994686754515833),
  @(0.6122060924509753, 0.2295189391074501, 0.05551053855075594, 0.1062326869378438, 2.866705290557533, 0.07973214163530429, 1.502624559269769, 0.1799087227181673, 0.8360799160353167, 0.3392015573179208, 2.999232255595508),
  @(0.6108539073057102, 0.2291870532635585, 0.05547243163406534, 0.1062094764981332, 2.866544158805482, 0.07973214163530429, 1.502705602965705, 0.1798966539493279, 0.8344314551528953, 0.3388107943481913, 2.999998630959212),
  @(0.6202736992587461, 0.2315035911224186, 0.05573924023191523, 0.106374682811559, 2.86777303936789, 0.07973214163530429, 1.502196356811154, 0.1799871166501106, 0.8459205195006234, 0.3415407656991789, 2.994747279165317),
  @(0.6626854175418373, 0.2420452178236872, 0.0569734139296898, 0.1072050326718959, 2.875929034661169, 0.07973214163530429, 1.501267957186954, 0.1805518433626361, 0.8977805665756193, 0.3540239374602052, 2.973260604696335),
  @(0.7488863686124887, 0.2638694768931202, 0.05959579168415274, 0.1091964447803804, 2.901751739211321, 0.07973214163530429, 1.504191089102669, 0.1822538565913092, 1.003655027482068, 0.3800704623262874, 2.937191387577698),
  @(0.8140741764138113, 0.2806132397455769, 0.06164426381100441, 0.1108802921018324, 2.9267284303114, 0.07973214163530429, 1.509237025052002, 0.1838669668882673, 1.084002982979285, 0.4001642874106892, 2.914393943654531),
  @(0.8441349944179137, 0.2883861370704892, 0.06260233349857458, 0.1116941493018189, 2.939398102018572, 0.07973214163530429, 1.51216346563784, 0.1846796235977237, 1.121115652017437, 0.4095140276668587, 2.904827861580642),
  @(0.8555767207447218, 0.2913520679209114, 0.06296886942995172, 0.1120092046634618, 2.944383975270298, 0.07973214163530429, 1.513362465987228, 0.1849986920389952, 1.13525018844561, 0.4130845771538958, 2.901321246922762),
  @(0.8531099470744437, 0.2907122998388729, 0.06288976373725319, 0.1119410468004567, 2.943301808790238, 0.07973214163530429, 1.513100199210683, 0.184929470976904, 1.132202471241754, 0.4123142615817059, 2.902071304036454),
  @(0.8450751434208428, 0.2886296946397522, 0.06263241395011221, 0.1117199316067534, 2.939804521514858, 0.07973214163530429, 1.512260287461203, 0.1847056464470711, 1.122276890207416, 0.4098071777994647, 2.904537047518744),
  @(0.8401611892975325, 0.2873569712158996, 0.0624752652975431, 0.1115853857046147, 2.937686837126506, 0.07973214163530429, 1.511757647259635, 0.1845700232117338, 1.116207703336698, 0.4082754229960059, 2.906062478721097),
  @(0.8121178022983031, 0.2801084083550336, 0.06158217626491336, 0.1108280662874463, 2.925926759590254, 0.07973214163530429, 1.509058479864883, 0.183815444001965, 1.08158888793227, 0.3995574610409705, 2.915035313083294),
  @(0.7950181387378734, 0.2757016667407299, 0.06104098427133664, 0.1103757223726269, 2.919047326553283, 0.07973214163530429, 1.507564304104697, 0.1833727246172003, 1.060495274545502, 0.394262762638462, 2.920746039301122),
  @(0.7852211383547569, 0.273181712445222, 0.06073217387625363, 0.1101200526046746, 2.915213534155967, 0.07973214163530429, 1.50676428320196, 0.1831255056457906, 1.0484156843942, 0.3912370627674235, 2.924106420233088),
  @(0.7819106194916685, 0.2723310199090747, 0.06062804103572716, 0.1100342617549934, 2.913936613159876, 0.07973214163530429, 1.506503607648838, 0.1830430763521136, 1.044334829292069, 0.3902159925305355, 2.925257188238874),
  @(0.7968344667633573, 0.2761692514380343, 0.06109833980055868, 0.1104234088802762, 2.919766915333142, 0.07973214163530429, 1.507717214491414, 0.1834190847895272, 1.062735251206504, 0.3948243567975993, 2.920130284660658),
  @(0.8474335774073154, 0.2892407950537006, 0.06270790269689286, 0.1117846923180679, 2.940826652537424, 0.07973214163530429, 1.512504524392661, 0.1847710816003598, 1.125190082055411, 0.4105427550616909, 2.903809653421064),
  @(0.8808430603184547, 0.2979149996755552, 0.06378160459404825, 0.1127143778317503, 2.955687105282422, 0.07973214163530429, 1.516162720722221, 0.1857207449623104, 1.166478722782699, 0.4209905147919315, 2.893818531908082),
  @(0.8629807143668984, 0.2932733889665542, 0.06320656957949922, 0.1122145322810937, 2.947655418055291, 0.07973214163530429, 1.514161803800746, 0.1852078491927003, 1.144399131689909, 0.4153983611092755, 2.8990891345628),
  @(0.7960131994960307, 0.2759578142250234, 0.0610724021122806, 0.1104018361385926, 2.919441211440727, 0.07973214163530429, 1.507647899917032, 0.1833981025973941, 1.061722410081103, 0.3945704030956279, 2.9204084268207),
  @(0.7252413710147039, 0.2578413610759185, 0.05886487493758352, 0.1086189053800446, 2.893712463016072, 0.07973214163530429, 1.502891803042679, 0.1817297464406664, 0.9745644049953057, 0.3728562075692992, 2.946299194559444)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $cellRef = $cols[$c] + $excelRow
        $ws.Range($cellRef).Value = $rowVals[$c]
    }
}
